$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "62.972.96"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -4.96%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.070.38"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -5.33%  "

$ws.Range("E4").Value = "  +0.15%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "547.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.60%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.31"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -12.53%  "

$ws.Range("E7").Value = "  +0.04%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.066.55"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.33%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.492"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -4.06%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -6.44%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.33"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -11.42%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.461"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -4.70%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "34.97"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -8.02%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000216"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -7.93%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.566.23"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.97%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "63.001.94"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -4.79%  "

$ws.Range("E17").Value = "  -2.85%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.068.57"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.83%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.66"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.36%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "485.58"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -10.84%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.45"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -7.17%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.705"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.27%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -7.06%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "77.79"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.40%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.26"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -8.95%  "

$ws.Range("E26").Value = "  +0.00%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.43"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -10.25%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -5.61%  "

$ws.Range("E29").Value = "  +0.14%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.95"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -13.84%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "26.30"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.80%  "

$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("E33").Value = "  -10.48%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "57.56"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.77%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "517.55"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -9.59%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.98"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.75%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.11"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -11.92%  "

$ws.Range("E38").Value = "  -13.30%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.095.65"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.60%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0795"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -8.17%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -5.74%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.09"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -6.04%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.64"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -14.12%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.253"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -7.59%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.07"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -12.63%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "120.36"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.62%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "24.52"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -8.41%  "

$ws.Range("E49").Value = "  -4.68%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.39"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +56.36%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0₃0498"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -10.75%  "
